$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date property value (row 8: A8="Date", B8=<timestamp>)
$ws.Range("B8").Value = "2025-07-11T12:29:53+00:00"

# Set the Jurisdiction value (row 11: A11="Jurisdiction", B11=<value>)
$ws.Range("B11").Value = "FRANCE"
